$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.553.59'
$ws.Range('E2').Value = '  -2.61%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.752.99'
$ws.Range('E3').Value = '  -3.42%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '324.10'
$ws.Range('E5').Value = '  -0.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4469'
$ws.Range('E7').Value = '  +2.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3624'
$ws.Range('E8').Value = '  -1.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07516'
$ws.Range('E9').Value = '  -1.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.15'
$ws.Range('E10').Value = '  -5.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.106'
$ws.Range('E11').Value = '  -2.94%  '
$ws.Range('E12').Value = '  +0.09%  '
$ws.Range('E13').Value = '  -5.75%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.047'
$ws.Range('E14').Value = '  -4.00%  '
$ws.Range('E15').Value = '  -4.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.751.66'
$ws.Range('E16').Value = '  -4.21%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '92.90'
$ws.Range('E17').Value = '  -2.30%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001064'
$ws.Range('E18').Value = '  -1.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06430'
$ws.Range('E19').Value = '  -1.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('E20').Value = '  +0.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.08'
$ws.Range('E21').Value = '  -1.80%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.848'
$ws.Range('E22').Value = '  -6.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.597.03'
$ws.Range('E23').Value = '  -2.51%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.25'
$ws.Range('E24').Value = '  -2.66%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.101'
$ws.Range('E25').Value = '  -0.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '162.53'
$ws.Range('E26').Value = '  +0.75%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.42'
$ws.Range('E27').Value = '  -1.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.952.25'
$ws.Range('E28').Value = '  -3.94%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.131'
$ws.Range('E29').Value = '  -6.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.60'
$ws.Range('E30').Value = '  -2.78%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.085'
$ws.Range('E31').Value = '  -10.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09016'
$ws.Range('E32').Value = '  -1.38%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.640'
$ws.Range('E33').Value = '  +1.86%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.550'
$ws.Range('E34').Value = '  -7.66%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.07'
$ws.Range('E35').Value = '  -6.75%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02332'
$ws.Range('E36').Value = '  -1.28%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2097'
$ws.Range('E37').Value = '  -3.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.6373'
$ws.Range('E38').Value = '  -3.26%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05975'
$ws.Range('E39').Value = '  -3.80%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.958'
$ws.Range('E40').Value = '  -5.09%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.207'
$ws.Range('E41').Value = '  +1.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9998'
$ws.Range('E42').Value = '  +0.14%  '
$ws.Range('E43').Value = '  -3.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.816'
$ws.Range('E44').Value = '  -3.22%  '
$ws.Range('E45').Value = '  -4.06%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5891'
$ws.Range('E46').Value = '  -3.43%  '
$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.714'
$ws.Range('E47').Value = '  -0.55%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.958'
$ws.Range('E48').Value = '  -2.79%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '121.50'
$ws.Range('E49').Value = '  -3.11%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.161'
$ws.Range('E50').Value = '  +0.39%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06867'
$ws.Range('E51').Value = '  -1.80%  '
